# Updated testcases for failing results and all URL
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "ProjectApproval" (sheet5.xml) -- edit this sheet first so that the
# workbook's active/selected tab ends up back on "TestCases" afterwards.
# ---------------------------------------------------------------------------
$wsPA = $wb.Worksheets.Item("ProjectApproval")
$wsPA.Range("G4").Value = "geoamps sample project"
$wsPA.Range("H4").Value = "geoamps sample project-00005"

$wsPA.Range("H4").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "TestCases" (sheet1.xml)
# ---------------------------------------------------------------------------
$wsTC = $wb.Worksheets.Item("TestCases")

$urlAlt = "https://uat.geoamps.com/altamps/Login/Login.aspx"
$urlDot = "https://uat.geoamps.com/dotAMPS/Login/Login.aspx"
$urlRow = "https://uat.geoamps.com/pipe/Login/Login.aspx"

# Rows 3 & 4 used to reference the old "environment1" shared string; they now
# point at the same "environment" label used by row 2.
$wsTC.Range("A3").Value = "environment"
$wsTC.Range("A4").Value = "environment"

# New rows 5-9: additional environment-specific test cases.
$wsTC.Range("A5").Value = "environmentALT"
$wsTC.Range("B5").Value = $urlAlt
$wsTC.Range("C5").Value = "ALT"

$wsTC.Range("A6").Value = "environmentDOT"
$wsTC.Range("B6").Value = $urlDot
$wsTC.Range("C6").Value = "DOT"

$wsTC.Range("A7").Value = "environmentROW"
$wsTC.Range("B7").Value = $urlRow
$wsTC.Range("C7").Value = "ROW"

$wsTC.Range("A8").Value = "environmentALTROW"
$wsTC.Range("B8").Value = $urlAlt
$wsTC.Range("C8").Value = "ALT"

$wsTC.Range("A9").Value = "environmentALTROW"
$wsTC.Range("B9").Value = $urlRow
$wsTC.Range("C9").Value = "ROW"

# Copy the look of the existing rows onto the new ones (style + number format).
$wsTC.Range("A2:C4").Copy() | Out-Null
$wsTC.Range("A5:C7").PasteSpecial(-4122) | Out-Null
$wsTC.Range("A2:C3").Copy() | Out-Null
$wsTC.Range("A8:C9").PasteSpecial(-4122) | Out-Null

# Hyperlinks for the new URL cells.
$wsTC.Hyperlinks.Add($wsTC.Range("B5"), $urlAlt) | Out-Null
$wsTC.Hyperlinks.Add($wsTC.Range("B6"), $urlDot) | Out-Null
$wsTC.Hyperlinks.Add($wsTC.Range("B7"), $urlRow) | Out-Null
$wsTC.Hyperlinks.Add($wsTC.Range("B8"), $urlAlt) | Out-Null
$wsTC.Hyperlinks.Add($wsTC.Range("B9"), $urlRow) | Out-Null

# Adding a hyperlink re-applies the built-in "Hyperlink" cell style; restore
# the plain style/value used by the rest of the column.
$wsTC.Range("B5").Value = $urlAlt
$wsTC.Range("B6").Value = $urlDot
$wsTC.Range("B7").Value = $urlRow
$wsTC.Range("B8").Value = $urlAlt
$wsTC.Range("B9").Value = $urlRow
$wsTC.Range("B2").Copy() | Out-Null
$wsTC.Range("B5:B9").PasteSpecial(-4122) | Out-Null

$wsTC.Range("A2").Select() | Out-Null
